$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update StartTime (column C) date serial values to reflect the newly
# added event file's schedule.
$ws.Range("C2").Value = 44059
$ws.Range("C3").Value = 44060
$ws.Range("C4").Value = 44061
$ws.Range("C5").Value = 44070
$ws.Range("C6").Value = 44071
$ws.Range("C7").Value = 44072
$ws.Range("C8").Value = 44073
$ws.Range("C9").Value = 44074
$ws.Range("C10").Value = 44075
$ws.Range("C14").Value = 44067
$ws.Range("C20").Value = 44068
